# Applies the "treasure hunt" edit:
#  - renames Sheet1 to "Wales"
#  - adds a new "North Ireland" sheet with treasure-hunt data
#  - moves the active/selected-tab state onto the new sheet

$wb = $excel.ActiveWorkbook

# --- Rename existing sheet -------------------------------------------------
$wsWales = $wb.Worksheets.Item(1)
$wsWales.Name = "Wales"
$wsWales.Range("B1").Select()

# --- Add the new sheet after Wales -----------------------------------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsWales)
$ws.Name = "North Ireland"

# --- Header row ---------------------------------------------------------
$headers = @("Location Coordinates", "Coordinates (Approximate)", "Treasure Value", "Likelihood %", "Recommend Reason", "Supporting Evidence")
for ($c = 1; $c -le 6; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = $headers[$c - 1]
}
$headerRange = $ws.Range("A1:F1")
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.WrapText = $true
$headerRange.Font.Name = "Segoe UI"
$headerRange.Font.Size = 9.6
$ws.Rows.Item(1).RowHeight = 57

# --- Data rows ---------------------------------------------------------
$data = @(
    @{ A = "Lacada Point, Portballintrae, Co. Antrim"; B = [char]0x00B0; C = "Very High"; D = 0.8;
       E = "Confirmed shipwreck location of Spanish Armada galleon Girona with documented treasure recovery";
       F = "Gold and artifacts from the Girona were discovered by Belgian diver Robert Sténuit in 19674";
       Link = "https://www.bbc.co.uk/news/uk-northern-ireland-40045580" },
    @{ A = "Near Kinbane Castle, Co. Antrim";
       C = "Medium"; D = 0.65;
       E = "Significant coin hoard suggests potential for additional caches during wartime period";
       F = "84 silver coins dating from 1551-1649 discovered in 2019, including rare Bavarian thaler5";
       Link = "https://www.bbc.co.uk/news/uk-northern-ireland-61914872" },
    @{ A = "Carrowdressex, Co. Down";
       C = "Medium"; D = 0.6;
       E = "Proximity to two 14th-century church sites with proven medieval artifacts";
       F = "Medieval gold ring brooch found in 2015 near historical church locations"; Sup = "36" },
    @{ A = "Ballyhornan Beach, Co. Down";
       C = "Medium"; D = 0.55;
       E = "Coastal location with confirmed 16th-century coin deposits";
       F = "Stack of twelve groats issued between 1555-1558 found in 2014"; Sup = "38" },
    @{ A = "Mara Castle area, Downpatrick, Co. Down";
       C = "Medium-High"; D = 0.7;
       E = "Viking settlement area with confirmed artifacts; arm rings typically found in hoards";
       F = "9th Century Viking arm ring discovered in pieces on farmland3";
       Link = "https://www.bbc.co.uk/news/uk-northern-ireland-37049242" },
    @{ A = "County Fermanagh";
       C = "High"; D = 0.6;
       E = "Evidence of Bronze Age wealth suggests potential for additional high-value deposits";
       F = "Bronze Age gold torc (720g) dating to 1,300 BC discovered in the" }
)

# Coordinates column (B) -- built separately to keep degree/quote characters exact
$coords = @(
    "55.2415$([char]0x00B0) N, 6.5167$([char]0x00B0) W",
    "55.2310$([char]0x00B0) N, 6.3360$([char]0x00B0) W",
    "54$([char]0x00B0) 16' 25`" N, 5$([char]0x00B0) 40' 36`" W",
    "54$([char]0x00B0)18'06.6875`" N, 5$([char]0x00B0)33'10.1583`" W",
    "54.32$([char]0x00B0) N, 5.72$([char]0x00B0) W",
    "54.34$([char]0x00B0) N, 7.63$([char]0x00B0) W"
)

$r = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $coords[$i]
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 4).NumberFormat = "0%"
    $ws.Cells.Item($r, 5).Value = $row.E

    $fCell = $ws.Cells.Item($r, 6)
    if ($row.Link) {
        $fCell.Value = $row.F
        $ws.Hyperlinks.Add($fCell, $row.Link, "", "", $row.F) | Out-Null
    } elseif ($row.Sup) {
        $fCell.Value = $row.F + $row.Sup
        $chars = $fCell.Characters($row.F.Length + 1, $row.Sup.Length)
        $chars.Font.Name = "Courier New"
        $chars.Font.Size = 9.6
    } else {
        $fCell.Value = $row.F
    }

    $dataRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 6))
    $dataRange.VerticalAlignment = -4108
    $dataRange.WrapText = $true
    $dataRange.Font.Name = "Segoe UI"
    $dataRange.Font.Size = 9.6

    if ($r -lt 7) {
        $dataRange.Borders.Item(9).LineStyle = 1
        $dataRange.Borders.Item(9).Weight = -4138
        $dataRange.Borders.Item(9).Color = 0
    }

    $r++
}

$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null

# --- View state: put the selection/active tab on the new sheet -------------
$ws.Range("A4").Select()
$ws.Activate()
